# Fruta / hortaliza, semanal
# Insert two new weekly price-report rows at the top of the "Vega Monumental
# Concepción - Palta" data block (row 191), pushing all the existing rows
# down by two (old row 191 -> new row 193, ..., old row 293 -> new row 295).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 191; this shifts every
# existing row (191-293) down to (193-295), matching the new dimension
# A1:T295 from the diff.
$ws.Rows("191:192").Insert()

# --- New row 191: Hass / Primera, new weekly observation (2021-09-21) ---
$ws.Range("A191").Value = 11
$ws.Range("B191").Value = "Vega Monumental Concepción"
$ws.Range("C191").Value = "Bíobío"
$ws.Range("D191").Value = 44460
$ws.Range("E191").Value = 8
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100106
$ws.Range("H191").Value = "Oleaginosos"
$ws.Range("I191").Value = 100106002
$ws.Range("J191").Value = "Palta"
$ws.Range("K191").Value = "Hass"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 50
$ws.Range("N191").Value = 3100
$ws.Range("O191").Value = 3100
$ws.Range("P191").Value = 3100
$ws.Range("Q191").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R191").Value = "Perú"
$ws.Range("S191").Value = 3100
$ws.Range("T191").Value = 1

# --- New row 192: Hass / Segunda, new weekly observation (2021-09-21) ---
$ws.Range("A192").Value = 11
$ws.Range("B192").Value = "Vega Monumental Concepción"
$ws.Range("C192").Value = "Bíobío"
$ws.Range("D192").Value = 44460
$ws.Range("E192").Value = 8
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100106
$ws.Range("H192").Value = "Oleaginosos"
$ws.Range("I192").Value = 100106002
$ws.Range("J192").Value = "Palta"
$ws.Range("K192").Value = "Hass"
$ws.Range("L192").Value = "Segunda"
$ws.Range("M192").Value = 50
$ws.Range("N192").Value = 2900
$ws.Range("O192").Value = 2900
$ws.Range("P192").Value = 2900
$ws.Range("Q192").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R192").Value = "Perú"
$ws.Range("S192").Value = 2900
$ws.Range("T192").Value = 1
